$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.607.58'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.25'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5311'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3971'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07757'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.07%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.119'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.43%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.00'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.18'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.325'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.595'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.001'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.824.60'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.09'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001088'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06608'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.80'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.092'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.616.54'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.237'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +6.47%  '
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.93'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.039.64'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.419'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.74'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.94%  '
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.733'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.21%  '
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07325'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.76%  '
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02353'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.903'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.27%  '
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("E40").Value = '  +2.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6285'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.396'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.56'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5931'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.19%  '
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.48'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("E49").Value = '  +4.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.193'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06955'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.90%  '
